# posSync packet fix: send x, y, z instead of just x, z
# (i.e. insert a posZ field right after the existing posY field in both
#  the CS_POSITION_SYNC and SC_POSITION_SYNC message definitions on the
#  "Player" sheet)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Player sheet - insert the new posZ rows
# ---------------------------------------------------------------
$wsPlayer = $wb.Worksheets.Item("Player")
$wsPlayer.Activate()

# CS_POSITION_SYNC block:
#   row9  = posX
#   row10 = posY
#   row11 = cameraYaw   <- insert new posZ row before this, pushing it to row12
$wsPlayer.Rows.Item(11).Insert()
$wsPlayer.Range("A11").Value = $wsPlayer.Range("A10").Text
$wsPlayer.Range("B11").Value = "posZ"
$wsPlayer.Range("C11").Value = $wsPlayer.Range("C10").Text

# SC_POSITION_SYNC block (shifted down by one row after the insert above):
#   row14 = playerID
#   row15 = posX
#   row16 = posY
#   row17 = cameraYaw   <- insert new posZ row before this, pushing it to row18
$wsPlayer.Rows.Item(17).Insert()
$wsPlayer.Range("A17").Value = $wsPlayer.Range("A16").Text
$wsPlayer.Range("B17").Value = "posZ"
$wsPlayer.Range("C17").Value = $wsPlayer.Range("C16").Text

# restore the view's selection on this sheet
$wsPlayer.Range("D24").Select()

# ---------------------------------------------------------------
# Room sheet - becomes the active tab / sheet when the workbook is reopened
# ---------------------------------------------------------------
$wsRoom = $wb.Worksheets.Item("Room")
$wsRoom.Activate()
$wsRoom.Range("A49").Select()
